$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column G with text values "1", "t2", "3" for rows 1-3,
# matching the data pattern added in the commit (input data fixes).
$ws.Range("G1").Value = "1"
$ws.Range("G2").Value = "t2"
$ws.Range("G3").Value = "3"
